$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.020.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.16%  '

$ws.Range("D3").Value = "'2.305.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = "'303.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("D6").Value = "'97.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("E7").Value = '  -1.94%  '

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").Value = "'0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.92%  '

$ws.Range("D10").Value = "'35.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.37%  '

$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("D12").Value = "'18.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.79%  '

$ws.Range("E13").Value = '  +1.34%  '

$ws.Range("D14").Value = "'6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.71%  '

$ws.Range("D15").Value = "'2.664.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("D16").Value = "'2.311.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").Value = "'42.865.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("D19").Value = "'12.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.40%  '

$ws.Range("D20").Value = "'0.0₃0897"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '

$ws.Range("E21").Value = '  -0.22%  '

$ws.Range("D22").Value = "'67.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").Value = "'236.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.64%  '

$ws.Range("D24").Value = "'2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.56%  '

$ws.Range("E25").Value = '  +0.16%  '

$ws.Range("D26").Value = "'2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").Value = "'24.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.46%  '

$ws.Range("E28").Value = '  +17.87%  '

$ws.Range("D29").Value = "'165.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("D31").Value = "'32.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.81%  '

$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("D33").Value = "'18.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.54%  '

$ws.Range("E34").Value = '  -1.01%  '

$ws.Range("E35").Value = '  -8.60%  '

$ws.Range("E36").Value = '  -1.45%  '

$ws.Range("D37").Value = "'0.0691"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.53%  '

$ws.Range("E38").Value = '  -0.52%  '

$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("E40").Value = '  +1.24%  '

$ws.Range("E41").Value = '  -0.68%  '

$ws.Range("D42").Value = "'1.997.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.98%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'10.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.90%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'0.0279"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.94%  '

$ws.Range("D45").Value = "'17.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.21%  '

$ws.Range("E46").Value = '  -1.47%  '

$ws.Range("D47").Value = "'2.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.08%  '

$ws.Range("D48").Value = "'2.531.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.26%  '

$ws.Range("E49").Value = '  -3.50%  '

$ws.Range("D50").Value = "'53.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.29%  '

$ws.Range("E51").Value = '  -0.50%  '
